$d = $word.ActiveDocument

$replacements = @(
    @{old="64×11="; new="92×40="},
    @{old="63×23="; new="14×73="},
    @{old="38×94="; new="34×16="},
    @{old="55×71="; new="49×21="},
    @{old="91×89="; new="71×68="},
    @{old="15×26="; new="39×81="},
    @{old="70×68="; new="74×82="},
    @{old="65×65="; new="85×98="},
    @{old="63×46="; new="16×66="},
    @{old="83×74="; new="34×81="},
    @{old="80×74="; new="79×73="},
    @{old="84×42="; new="62×64="},
    @{old="31×65="; new="53×37="},
    @{old="30×18="; new="47×45="},
    @{old="67×36="; new="62×21="},
    @{old="40×53="; new="45×11="},
    @{old="86×26="; new="59×56="},
    @{old="89×82="; new="66×95="},
    @{old="73×24="; new="80×46="},
    @{old="41×60="; new="78×21="},
    @{old="54×92="; new="98×18="},
    @{old="60×82="; new="76×12="},
    @{old="90×54="; new="94×80="},
    @{old="32×86="; new="69×91="},
    @{old="86×16="; new="98×15="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
